# PSP_Sheet_2조.xlsx — add 6 new Time Recording Log entries (rows 21-26)
# to the "작성자명" sheet (first worksheet), matching the target diff:
#   - rows 21-24: fill in date/start/stop/interruption/delta/activity
#   - rows 25-26: two new entries for 2019-11-30 (같은 날짜, 두 구간)
#   - small style touches: font on F22/F23, thinner side-borders on F25/F26,
#     no outer border on the new B25:C26 time cells
#   - move the saved view's selection to F26

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Row 21 — 2019-11-19 (serial 43788)
# ---------------------------------------------------------------------
$ws.Range("A21").Value = 43788
$ws.Range("B21").Value = 0.5625
$ws.Range("C21").Value = 0.61458333333333337
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 75
$ws.Range("F21").Value = "개발 환경 세팅"

# ---------------------------------------------------------------------
# Row 22 — 2019-11-21 (serial 43790)
# ---------------------------------------------------------------------
$ws.Range("A22").Value = 43790
$ws.Range("B22").Value = 0.5625
$ws.Range("C22").Value = 0.61458333333333337
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 75
$ws.Range("F22").Value = "최종 결과물을 위한 업무 분배"
$ws.Range("F22").Font.Name = "맑은 고딕"
$ws.Range("F22").Font.Size = 10

# ---------------------------------------------------------------------
# Row 23 — 2019-11-26 (serial 43795)
# ---------------------------------------------------------------------
$ws.Range("A23").Value = 43795
$ws.Range("B23").Value = 0.5625
$ws.Range("C23").Value = 0.61458333333333337
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 75
$ws.Range("F23").Value = "요람 데이터 정제 작업"
$ws.Range("F23").Font.Name = "맑은 고딕"
$ws.Range("F23").Font.Size = 10
$ws.Rows.Item(23).RowHeight = 16

# ---------------------------------------------------------------------
# Row 24 — 2019-11-28 (serial 43797)
# ---------------------------------------------------------------------
$ws.Range("A24").Value = 43797
$ws.Range("B24").Value = 0.5625
$ws.Range("C24").Value = 0.61458333333333337
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 75
$ws.Range("F24").Value = "개발, 문서작성 등 개별 업무 진행"
# "개발" stays in the cell's base font; the remainder switches to 돋움
$ws.Range("F24").Characters(3, 17).Font.Name = "돋움"
$ws.Range("F24").Characters(3, 17).Font.Size = 10

# ---------------------------------------------------------------------
# Row 25 — 2019-11-30 (serial 43799), first block of the day
# ---------------------------------------------------------------------
$ws.Range("A25").Value = 43799
$ws.Range("B25").Value = 0.41666666666666669
$ws.Range("C25").Value = 0.5
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 120
# B25/C25 keep the h:mm number format but lose the surrounding box border
$ws.Range("B25:C25").Borders.LineStyle = -4142

$f25text = "개별 업무 진행(각종 정보 입력 및 저장 구현, HTML/CSS 작업, 추천 로직 구현, 문서 작성)"
$ws.Range("F25").Value = $f25text
$ws.Range("F25").Characters(1, 27).Font.Name = "맑은 고딕"
$ws.Range("F25").Characters(1, 27).Font.Size = 10
$ws.Range("F25").Characters(28, 9).Font.Name = "Arial"
$ws.Range("F25").Characters(28, 9).Font.Size = 10
$ws.Range("F25").Characters(37, 20).Font.Name = "맑은 고딕"
$ws.Range("F25").Characters(37, 20).Font.Size = 10
# F25 border becomes thin left/right only (top/bottom cleared)
$ws.Range("F25").Borders(8).LineStyle = -4142
$ws.Range("F25").Borders(9).LineStyle = -4142
$ws.Range("F25").Borders(7).LineStyle = 1
$ws.Range("F25").Borders(10).LineStyle = 1
$ws.Rows.Item(25).RowHeight = 16

# ---------------------------------------------------------------------
# Row 26 — 2019-11-30 (serial 43799), second block of the day
# ---------------------------------------------------------------------
$ws.Range("A26").Value = 43799
$ws.Range("B26").Value = 0.54166666666666663
$ws.Range("C26").Value = 0.91666666666666663
$ws.Range("D26").Value = 120
$ws.Range("E26").Value = 420
$ws.Range("B26:C26").Borders.LineStyle = -4142

$f26text = $f25text
$ws.Range("F26").Value = $f26text
$ws.Range("F26").Characters(1, 27).Font.Name = "맑은 고딕"
$ws.Range("F26").Characters(1, 27).Font.Size = 10
$ws.Range("F26").Characters(28, 9).Font.Name = "Arial"
$ws.Range("F26").Characters(28, 9).Font.Size = 10
$ws.Range("F26").Characters(37, 20).Font.Name = "맑은 고딕"
$ws.Range("F26").Characters(37, 20).Font.Size = 10
$ws.Range("F26").Borders(8).LineStyle = -4142
$ws.Range("F26").Borders(9).LineStyle = -4142
$ws.Range("F26").Borders(7).LineStyle = 1
$ws.Range("F26").Borders(10).LineStyle = 1
$ws.Rows.Item(26).RowHeight = 16

# ---------------------------------------------------------------------
# View: scroll so row 14 is at the top and select F26 (last edited cell)
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 14
$ws.Range("F26").Select()
